# Andy altered income ratio slide
$p = $ppt.ActivePresentation

# --- Slide 14 ("Future Work"): merge the two trailing runs of the
#     "Include weighting..." bullet into a single run ---
$s14 = $p.Slides.Item(14)
$contentShape = $s14.Shapes.Item("Content Placeholder 2")
$para = $contentShape.TextFrame.TextRange.Paragraphs(7)
$run2 = $para.Runs(2)
$run3 = $para.Runs(3)
$run2.Text = $run2.Text + $run3.Text
$run3.Text = ""

# --- Slide 8 ("Findings"): resize/reposition the picture and remove the
#     "Content Placeholder 2" bullet list shape entirely ---
$s8 = $p.Slides.Item(8)
$pic = $s8.Shapes.Item("Picture 2")
$pic.Left = 90
$pic.Top = 12
$pic.Width = 552
$pic.Height = 551.1695

$bullets = $s8.Shapes.Item("Content Placeholder 2")
$bullets.Cut()
